$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'272.74"
$ws.Range("E2").Value = "'0.05%"
$ws.Range("D3").Value = "'26.80"
$ws.Range("E3").Value = "'-0.54%"
$ws.Range("D4").Value = "'4.909"
$ws.Range("E4").Value = "'3.96%"
$ws.Range("D5").Value = "'0.06321"
$ws.Range("E5").Value = "'2.87%"
$ws.Range("D6").Value = "'6.906"
$ws.Range("E6").Value = "'2.60%"
$ws.Range("D7").Value = "'3.359"
$ws.Range("E7").Value = "'5.54%"
$ws.Range("D8").Value = "'1.280"
$ws.Range("E8").Value = "'40.97%"
$ws.Range("D9").Value = "'0.8869"
$ws.Range("E9").Value = "'2.76%"
$ws.Range("D10").Value = "'0.1461"
$ws.Range("E10").Value = "'1.27%"
$ws.Range("E11").Value = "'0.86%"
$ws.Range("D12").Value = "'0.07404"
$ws.Range("E12").Value = "'3.51%"
$ws.Range("D13").Value = "'0.03138"
$ws.Range("E13").Value = "'-1.21%"
$ws.Range("D14").Value = "'0.09036"
$ws.Range("E14").Value = "'-0.01%"
$ws.Range("D15").Value = "'0.001544"
$ws.Range("E15").Value = "'0.29%"
$ws.Range("D16").Value = "'0.0006326"
$ws.Range("E16").Value = "'4.45%"
$ws.Range("D17").Value = "'0.006042"
$ws.Range("E17").Value = "'2.15%"
$ws.Range("D18").Value = "'3.462"
$ws.Range("E18").Value = "'0.03%"
$ws.Range("D19").Value = "'2.272"
$ws.Range("E19").Value = "'0.32%"
$ws.Range("E20").Value = "'2.63%"
$ws.Range("D21").Value = "'0.1336"
$ws.Range("E21").Value = "'3.09%"
$ws.Range("D22").Value = "'3.904"
$ws.Range("E22").Value = "'1.61%"
$ws.Range("D23").Value = "'0.04360"
$ws.Range("E23").Value = "'3.04%"
$ws.Range("D24").Value = "'0.001178"
$ws.Range("E24").Value = "'0.20%"
$ws.Range("D25").Value = "'0.003672"
$ws.Range("E25").Value = "'-12.54%"
$ws.Range("E26").Value = "'0.10%"
$ws.Range("D27").Value = "'0.0001700"
$ws.Range("E27").Value = "'1.72%"
$ws.Range("D40").Value = "'0.04027"
$ws.Range("E40").Value = "'1.47%"
$ws.Range("D41").Value = "'0.006628"
$ws.Range("E41").Value = "'6.54%"
$ws.Range("E42").Value = "'2.81%"
$ws.Range("E43").Value = "'-3.14%"
$ws.Range("E44").Value = "'-5.93%"
$ws.Range("D45").Value = "'0.00005318"
$ws.Range("E45").Value = "'3.77%"
$ws.Range("D46").Value = "'2.365"
$ws.Range("E46").Value = "'163.64%"
$ws.Range("D47").Value = "'0.02605"
$ws.Range("E47").Value = "'-12.66%"
